$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update patient name and email test data
$ws.Range("A2").Value = "19Dec"
$ws.Range("B2").Value = "TestPatient"
$ws.Range("C2").Value = "Dec@19.com"

# Update the hyperlink target/display text to match new email
$ws.Hyperlinks.Item(1).Address = "mailto:Dec@19.com"
$ws.Hyperlinks.Item(1).TextToDisplay = "Dec@19.com"

# Update the active cell selection
$ws.Range("D2").Select()
